# "Added change password scenarios 2"
# Adds a new "changepassword" worksheet (after the existing "forgotpassword"
# sheet) containing new change-password test-case scenarios, and makes it
# the active sheet/tab.

$wb = $excel.ActiveWorkbook

# --- Create the new worksheet after the last existing sheet -----------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "changepassword"

# --- Column A-C: test case name / input values -------------------------------
# (filled in left-to-right, top-to-bottom first, matching how the data was
# originally authored)
$ws.Range("A1").Value = "*** Test Cases ***"
$ws.Range("B1").Value = '${newpassword}'
$ws.Range("C1").Value = '${confirmpassword}'

$ws.Range("A2").Value = "Blank new password"
$ws.Range("C2").Value = "Tester101!"

$ws.Range("A3").Value = "Less than 8 characters"
$ws.Range("B3").Value = "Test1!"
$ws.Range("C3").Value = "Test1!"

$ws.Range("A4").Value = "Password did not meet the requirements"
$ws.Range("B4").Value = "Tester101"
$ws.Range("C4").Value = "Tester101"

$ws.Range("A5").Value = "Password mismatch"
$ws.Range("B5").Value = "Tester101!"
$ws.Range("C5").Value = "Tester101!!"

# --- Column D: expected error message (filled in afterwards) ----------------
$ws.Range("D1").Value = '${errormessage}'
$ws.Range("D2").Value = "This field is required"
$ws.Range("D3").Value = "Minimum of 8 digits"
$ws.Range("D4").Value = "It must contain at least one lowercase letter, one uppercase letter, one numeric digit, and one special character"
$ws.Range("D5").Value = "Password do not match"

# D2's message wraps onto multiple lines
$ws.Range("D2").WrapText = $true

# --- Column widths ------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 38.29
$ws.Columns.Item(2).ColumnWidth = 37.29
$ws.Columns.Item(3).ColumnWidth = 36
$ws.Columns.Item(4).ColumnWidth = 102.14

# --- Sheet view / selection ----------------------------------------------------
$ws.Range("C18").Select()

# --- Page setup ------------------------------------------------------------
$ws.PageSetup.Orientation = 1

# --- forgotpassword sheet selection updates when it is no longer the active tab
$wsForgot = $wb.Worksheets.Item("forgotpassword")
$wsForgot.Range("A1:C1").Select()

# Make the new sheet the active / visible tab
$ws.Select()
